# Fixed naive component forecaster bug - Presentation state 11.02.
#
# The rolling error-statistics table (columns B:G, rows 2-11) is shifted
# down by one row (oldest quarter's row, row 11, drops off) and a newly
# computed row of statistics is inserted at row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 2-10 down into rows 3-11 (columns B:G), working
# from the bottom up so we don't overwrite data before reading it.
for ($r = 10; $r -ge 2; $r--) {
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($r + 1, $c).Value2 = $ws.Cells.Item($r, $c).Value2
    }
}

# Insert the newly computed statistics into row 2.
$ws.Range("B2").Value2 = 0.1279039517177195
$ws.Range("C2").Value2 = 1.47408018424297
$ws.Range("D2").Value2 = 8.675642101351979
$ws.Range("E2").Value2 = 2.945444296087091
$ws.Range("F2").Value2 = 3.008801495394956
$ws.Range("G2").Value2 = 23
